$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target table: all Sending-cluster x Target-cluster combinations for
# ECs / FAPs / MuSCs with Ligand Bmp15 / Receptor Bmpr1a (9 rows total),
# replacing the old 3-row table that always used "FAPs" as sending cluster.
$rows = @(
    @("ECs","Bmp15","Bmpr1a","ECs",1,0.3333333333333333,0.02770266666666667,0.083108,0.05180130905700151,0.05180130905700151,3,1,1.443038,4.329114,0.0289666880885598,0.0289666880885598,0.03997600070133333,0.359784006312,0.00150051236203325,0.00150051236203325),
    @("ECs","Bmp15","Bmpr1a","FAPs",1,0.3333333333333333,0.02770266666666667,0.083108,0.05180130905700151,0.05180130905700151,3,1,29.20351433333333,87.610543,0.5862140087672342,0.5862140087672342,0.8090152230715556,7.281137007644,0.03036665304169529,0.03036665304169529),
    @("ECs","Bmp15","Bmpr1a","MuSCs",1,0.3333333333333333,0.02770266666666667,0.083108,0.05180130905700151,0.05180130905700151,3,1,19.170603,57.511809,0.384819303144206,0.384819303144206,0.531076824708,4.779691422372,0.01993414365327297,0.01993414365327297),
    @("FAPs","Bmp15","Bmpr1a","ECs",1,0.3333333333333333,0.4126986666666667,1.238096,0.7717066171516261,0.7717066171516261,3,1,1.443038,4.329114,0.0289666880885598,0.0289666880885598,0.5955398585493333,5.359858726944,0.02235378487490879,0.02235378487490878),
    @("FAPs","Bmp15","Bmpr1a","FAPs",1,0.3333333333333333,0.4126986666666667,1.238096,0.7717066171516261,0.7717066171516261,3,1,29.20351433333333,87.610543,0.5862140087672342,0.5862140087672342,12.05225142734756,108.470262846128,0.452385229632656,0.452385229632656),
    @("FAPs","Bmp15","Bmpr1a","MuSCs",1,0.3333333333333333,0.4126986666666667,1.238096,0.7717066171516261,0.7717066171516261,3,1,19.170603,57.511809,0.384819303144206,0.384819303144206,7.911682297296001,71.205140675664,0.2969676026440614,0.2969676026440613),
    @("MuSCs","Bmp15","Bmpr1a","ECs",2,0.6666666666666666,0.09438566666666666,0.283157,0.1764920737913724,0.1764920737913724,3,1,1.443038,4.329114,0.0289666880885598,0.0289666880885598,0.1362021036553333,1.225818932898,0.005112390851617764,0.005112390851617763),
    @("MuSCs","Bmp15","Bmpr1a","FAPs",2,0.6666666666666666,0.09438566666666666,0.283157,0.1764920737913724,0.1764920737913724,3,1,29.20351433333333,87.610543,0.5862140087672342,0.5862140087672342,2.756393169361222,24.807538524251,0.1034621260928829,0.1034621260928829),
    @("MuSCs","Bmp15","Bmpr1a","MuSCs",2,0.6666666666666666,0.09438566666666666,0.283157,0.1764920737913724,0.1764920737913724,3,1,19.170603,57.511809,0.384819303144206,0.384819303144206,1.809430144557,16.284871301013,0.06791755684687172,0.0679175568468717)
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $rowData = $rows[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $addr = $cols[$j] + $r
        $ws.Range($addr).Value = $rowData[$j]
    }
}

Write-Host "updated $($rows.Length) rows"
